$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# --- Project title ---
$ws.Range("B1").Value = "Complexity Analyzer"

# --- Project lead name (B3) and Project Start date (E3, now a fixed date instead of =TODAY()) ---
$ws.Range("B3").Value = "Aaron Downing "
$ws.Range("E3").Value = 45939

# --- New cell B4: second team-member name, using a 14pt font like the other label cells ---
$ws.Range("B4").Value = "Ryerson Brower "
$ws.Range("B4").Font.Size = 14

# --- Task 1 (row 9): renamed "Task 1" -> "Coding", assigned to "All", 100% progress,
#     and duration extended from 3 to 10 days ---
$ws.Range("B9").Value = "Coding"
$ws.Range("C9").Value = "All"
$ws.Range("D9").Value = 1
$ws.Range("F9").Formula = "=E9+10"

# --- Task 2 (row 10): renamed "Task 2" -> "Tech Doc", assigned to "All", 100% progress,
#     and duration extended from 2 to 7 days ---
$ws.Range("B10").Value = "Tech Doc"
$ws.Range("C10").Value = "All"
$ws.Range("D10").Value = 1
$ws.Range("F10").Formula = "=E10+7"

# --- Task 3 (row 11): renamed "Task 3" -> "Org Doc", assigned to "All", 100% progress,
#     start/end dates chained off the previous task (row 10) ---
$ws.Range("B11").Value = "Org Doc"
$ws.Range("C11").Value = "All"
$ws.Range("D11").Value = 1
$ws.Range("E11").Formula = "=E10"
$ws.Range("F11").Formula = "=F10"

# --- Task 4 (row 12): renamed "Task 4" -> "Presntation ", assigned to "All", 100% progress,
#     start/end dates chained off the previous task (row 11) ---
$ws.Range("B12").Value = "Presntation "
$ws.Range("C12").Value = "All"
$ws.Range("D12").Value = 1
$ws.Range("E12").Formula = "=E11"
$ws.Range("F12").Formula = "=F11"

# --- Update the active selection to match the saved view state ---
$ws.Range("C3:D3").Select()
